# Append 5 new match rows (22-26) to Sheet1, mirroring the style of the
# last existing data row (row 21: bold/bordered index col, date-formatted
# kickoff col).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing last row as a formatting template for the new rows so
# the "Indice" column (bold, centered, bordered) and the "data_partida"
# column (yyyy-mm-dd hh:mm:ss number format) keep matching styles.
$ws.Range("A21:V21").Copy($ws.Range("A22:V26"))

$newRows = @(
    @{ Row=22; Indice=21; Data=45192.70833333334; Home="Hapoel Petah Tikva"; HomeGols=1; Away="H. Beer Sheva"; AwayGols=0;
       HOdds=3.35; HDate="21/09/2023 03:42"; HCOdds=6.79; HCDate="23/09/2023 16:59";
       DOdds=3.33; DDate="21/09/2023 03:42"; DCOdds=4.21; DCDate="23/09/2023 16:59";
       AOdds=2.2;  ADate="21/09/2023 03:42"; ACOdds=1.5; ACDate="23/09/2023 16:52";
       Url="https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-petah-tikva-h-beer-sheva/8jktLzDi/" },

    @{ Row=23; Indice=22; Data=45192.77083333334; Home="Hapoel Hadera"; HomeGols=1; Away="Hapoel Haifa"; AwayGols=2;
       HOdds=3.08; HDate="18/09/2023 18:12"; HCOdds=3.24; HCDate="23/09/2023 18:23";
       DOdds=3.22; DDate="18/09/2023 18:12"; DCOdds=3.4; DCDate="23/09/2023 18:28";
       AOdds=2.38; ADate="18/09/2023 18:12"; ACOdds=2.26; ACDate="23/09/2023 18:23";
       Url="https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-hadera-hapoel-haifa/l6YTH4oU/" },

    @{ Row=24; Indice=23; Data=45192.78125; Home="Hapoel Jerusalem"; HomeGols=1; Away="Netanya"; AwayGols=1;
       HOdds=2.7;  HDate="17/09/2023 18:12"; HCOdds=2.66; HCDate="23/09/2023 18:44";
       DOdds=3.21; DDate="17/09/2023 18:12"; DCOdds=3.22; DCDate="23/09/2023 18:37";
       AOdds=2.58; ADate="17/09/2023 18:12"; ACOdds=2.8; ACDate="23/09/2023 18:44";
       Url="https://www.betexplorer.com/football/israel/ligat-ha-al/hapoel-jerusalem-netanya/2DlpKfSc/" },

    @{ Row=25; Indice=24; Data=45192.78125; Home="SC Ashdod"; HomeGols=1; Away="Maccabi Petah Tikva"; AwayGols=1;
       HOdds=2.16; HDate="20/09/2023 15:12"; HCOdds=2.37; HCDate="23/09/2023 18:43";
       DOdds=3.31; DDate="20/09/2023 15:12"; DCOdds=3.32; DCDate="23/09/2023 18:42";
       AOdds=3.45; ADate="20/09/2023 15:12"; ACOdds=3.12; ACDate="23/09/2023 18:43";
       Url="https://www.betexplorer.com/football/israel/ligat-ha-al/sc-ashdod-maccabi-petah-tikva/fqZPIpWN/" },

    @{ Row=26; Indice=25; Data=45192.79166666666; Home="Maccabi Bnei Raina"; HomeGols=2; Away="Beitar Jerusalem"; AwayGols=1;
       HOdds=2.68; HDate="18/09/2023 18:12"; HCOdds=3.15; HCDate="23/09/2023 18:56";
       DOdds=3.25; DDate="18/09/2023 18:12"; DCOdds=3.34; DCDate="23/09/2023 18:56";
       AOdds=2.68; ADate="18/09/2023 18:12"; ACOdds=2.34; ACDate="23/09/2023 18:56";
       Url="https://www.betexplorer.com/football/israel/ligat-ha-al/maccabi-bnei-raina-beitar-jerusalem/jFKeDOgo/" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Indice
    $ws.Cells.Item($row, 2).Value = "israel"
    $ws.Cells.Item($row, 3).Value = "ligat-ha-al"
    $ws.Cells.Item($row, 4).Value = "2023-2024"
    $ws.Cells.Item($row, 5).Value = $r.Data
    $ws.Cells.Item($row, 6).Value = $r.Home
    $ws.Cells.Item($row, 7).Value = $r.HomeGols
    $ws.Cells.Item($row, 8).Value = $r.Away
    $ws.Cells.Item($row, 9).Value = $r.AwayGols

    $ws.Cells.Item($row, 10).Value = $r.HOdds
    $ws.Cells.Item($row, 11).Value = $r.HDate
    $ws.Cells.Item($row, 12).Value = $r.HCOdds
    $ws.Cells.Item($row, 13).Value = $r.HCDate

    $ws.Cells.Item($row, 14).Value = $r.DOdds
    $ws.Cells.Item($row, 15).Value = $r.DDate
    $ws.Cells.Item($row, 16).Value = $r.DCOdds
    $ws.Cells.Item($row, 17).Value = $r.DCDate

    $ws.Cells.Item($row, 18).Value = $r.AOdds
    $ws.Cells.Item($row, 19).Value = $r.ADate
    $ws.Cells.Item($row, 20).Value = $r.ACOdds
    $ws.Cells.Item($row, 21).Value = $r.ACDate

    $ws.Cells.Item($row, 22).Value = $r.Url
}
